# Append newly logged model-evaluation runs (rows 18-29) to the "Metrics"
# sheet and tidy up the stray empty "Test R2" placeholder cell that used to
# mark the most-recent row (it now belongs on the new last row, 29).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{Row=18; A="2024-12-02 14:28:32"; B=0.9961272683622816;  C=0.007608786382250166;  D=0.0001768543921761174;  E=0.01329866129263083;  F="NONE";  G=0.01506293455905632;  H=0.0002917397529903427; I=0.01708039089102889},
    @{Row=19; A="2024-12-02 14:33:35"; B=0.9961272683622816;  C=0.007608786382250166;  D=0.0001768543921761174;  E=0.01329866129263083;  F="NONE";  G=0.01506293455905632;  H=0.0002917397529903427; I=0.01708039089102889},
    @{Row=20; A="2024-12-02 16:47:06"; B=0.9956057609533321;  C=0.008807628479416482;  D=0.0002013509695350387;  E=0.01418981922136567;  F="NONE";  G=0.01445152294439398;  H=0.000272907172521264;  I=0.01651990231573008},
    @{Row=21; A="2024-12-02 16:51:42"; B=0.9956057609533321;  C=0.008807628479416482;  D=0.0002013509695350387;  E=0.01418981922136567;  F="NONE";  G=0.01445152294439398;  H=0.000272907172521264;  I=0.01651990231573008},
    @{Row=22; A="2024-12-02 16:55:08"; B=0.9921798356098381;  C=0.0139442818635473;    D=0.0003583287359029029;  E=0.01892957305125773;  F="NONE";  G=0.01143688690466009;  H=0.0002063415543995868; I=0.01436459377774348},
    @{Row=23; A="2024-12-02 17:03:24"; B=0.9921798356098381;  C=0.0139442818635473;    D=0.0003583287359029029;  E=0.01892957305125773;  F="NONE";  G=0.01143688690466009;  H=0.0002063415543995868; I=0.01436459377774348},
    @{Row=24; A="2024-12-02 17:07:43"; B=0.996048209643299;   C=0.007774794003715868;  D=0.0001810774493781663;  E=0.01345650212269765;  F="NONE";  G=0.01322528850835882;  H=0.0002649595870469305; I=0.01627757927478563},
    @{Row=25; A="2024-12-02 17:10:21"; B=0.996048209643299;   C=0.007774794003715868;  D=0.0001810774493781663;  E=0.01345650212269765;  F="NONE";  G=0.01322528850835882;  H=0.0002649595870469305; I=0.01627757927478563},
    @{Row=26; A="2024-12-02 17:26:24"; B=0.996048209643299;   C=0.007774794003715868;  D=0.0001810774493781663;  E=0.01345650212269765;  F="NONE";  G=0.01322528850835882;  H=0.0002649595870469305; I=0.01627757927478563},
    @{Row=27; A="2024-12-02 17:32:08"; B=0.996048209643299;   C=0.007774794003715868;  D=0.0001810774493781663;  E=0.01345650212269765;  F="NONE";  G=0.01322528850835882;  H=0.0002649595870469305; I=0.01627757927478563},
    @{Row=28; A="2024-12-02 17:39:14"; B=0.9960491470277466;  C=0.007710636910772076;  D=0.0001810345700240405;  E=0.01345490877055807;  F="NONE";  G=0.01676170360845648;  H=0.0003620553488139641; I=0.01902775206938445},
    @{Row=29; A="2024-12-03 10:35:11"; B=0.995863692216895;   C=0.008396577694415567;  D=0.0001895323645330481;  E=0.01376707538052466;  F="EMPTY"; G=0.01337260734838577;  H=0.0002468636198265327; I=0.01571189421510127}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A   # Date
    $ws.Cells.Item($row, 2).Value = $r.B   # Train R2
    $ws.Cells.Item($row, 3).Value = $r.C   # Train MAE
    $ws.Cells.Item($row, 4).Value = $r.D   # Train MSE
    $ws.Cells.Item($row, 5).Value = $r.E   # Train RMSE
    if ($r.F -eq "EMPTY") {
        # Test R2 not yet available for this run: keep the cell present but
        # blank, matching how the sheet previously marked its newest row.
        $ws.Cells.Item(17, 6).Copy($ws.Cells.Item($row, 6))
    }
    $ws.Cells.Item($row, 7).Value = $r.G   # Test MAE
    $ws.Cells.Item($row, 8).Value = $r.H   # Test MSE
    $ws.Cells.Item($row, 9).Value = $r.I   # Test RMSE
}

# Row 17 is no longer the newest row, so drop its placeholder blank Test R2
# cell entirely (the column is simply absent on older rows).
$ws.Cells.Item(17, 6).ClearContents()
